$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate rows 2-5 with the new retailer data. The order of assignment
# below matches the order the values were originally entered so that the
# shared-strings table is rebuilt with the same ordering.
$ws.Range("B2").Value = "Mobile Palace "
$ws.Range("C2").Value = "School MarketPuthiaRajshahi"
$ws.Range("A2").Value = "RET-08822"

$ws.Range("B3").Value = "Ma Digital Studio & Electronics"
$ws.Range("A3").Value = "RET-35304"
$ws.Range("C3").Value = "Edilpur Gurudaspur Natore."

$ws.Range("A4").Value = "RET-26510"
$ws.Range("B4").Value = "Mahfuz Telecom "
$ws.Range("C4").Value = "Kakramari Bazar Charghat Rajshahi "

$ws.Range("B5").Value = "Abir Electronics"
$ws.Range("A5").Value = "RET-36274"
$ws.Range("C5").Value = "Kaligong Bazar Singra Natore"

# Rows 6-11: clear content (retain/restore formatting via borders)
$ws.Range("A6:C11").ClearContents()

# Update selection
$ws.Range("G10").Select()
